# Jesse_James_TE_2018 scraper update: add "height" and "weight" columns
# between the existing "fumbles" column and the "fantasy points" column.
#
# Layout before:  A=index B=rec_yds C=rec_td D=fumbles E=fantasy points
# Layout after:   A=index B=rec_yds C=rec_td D=fumbles E=height F=weight G=fantasy points

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at E:F. This shifts the existing "fantasy points"
# column (and all of its header/data/style) from E to G.
$ws.Range("E1:F1").EntireColumn.Insert()

# Header cells for the two newly inserted columns. EntireColumn.Insert()
# already carries over the bold/centered/bordered header style from the
# neighboring column, so only the text needs to be set here.
$ws.Range("E1").Value2 = "height"
$ws.Range("F1").Value2 = "weight"

# Data values for the new columns - every player row carries the same
# height/weight constant from the scrape.
$ws.Range("E2:E17").Value2 = 6.583333333333333
$ws.Range("F2:F17").Value2 = 250
